$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 3819.8572
$ws.Range("J29").Value = 6999
$ws.Range("L29").Value = 20997
$ws.Range("N29").Value = -21559
$ws.Range("H38").Value = 1527.5555
$ws.Range("I38").Value = 1249.75
$ws.Range("K38").Value = 3749.25
$ws.Range("M38").Value = -3377.25
$ws.Range("H43").Value = 25002876
$ws.Range("I43").Value = 50001750
$ws.Range("K43").Value = 50001750
$ws.Range("M43").Value = -50001681
$ws.Range("H53").Value = 239.375
$ws.Range("J53").Value = 85.5
$ws.Range("L53").Value = 85.5
$ws.Range("N53").Value = -1359.5
$ws.Range("H58").Value = 1673.909
$ws.Range("J58").Value = 2733.1667
$ws.Range("L58").Value = 8199.500100000001
$ws.Range("N58").Value = -8499.500100000001
$ws.Range("H64").Value = 6772.636
$ws.Range("I64").Value = 4499.8335
$ws.Range("J64").Value = 9500
$ws.Range("K64").Value = 4499.8335
$ws.Range("L64").Value = 9500
$ws.Range("M64").Value = -4251.8335
$ws.Range("N64").Value = -9996
$ws.Range("H67").Value = 6772.636
$ws.Range("I67").Value = 4499.8335
$ws.Range("J67").Value = 9500
$ws.Range("K67").Value = 4499.8335
$ws.Range("L67").Value = 9500
$ws.Range("M67").Value = -3641.8335
$ws.Range("N67").Value = -11216
$ws.Range("H103").Value = 2015.2858
$ws.Range("I103").Value = 3556.5
$ws.Range("J103").Value = 1398.8
$ws.Range("K103").Value = 10669.5
$ws.Range("L103").Value = 4196.4
$ws.Range("M103").Value = -10083.5
$ws.Range("N103").Value = -5368.4
$ws.Range("H113").Value = 4139
$ws.Range("I113").Value = 2712
$ws.Range("K113").Value = 2712
$ws.Range("M113").Value = 542
$ws.Range("H116").Value = 6665.6665
$ws.Range("I116").Value = 5665.1665
$ws.Range("K116").Value = 5665.1665
$ws.Range("M116").Value = -2223.1665
$ws.Range("H132").Value = 8172.5
$ws.Range("I132").Value = 12345
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 37035
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -34505
$ws.Range("N132").Value = -17060
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1178.2
$ws.Range("I2").Value = 1178.2
$ws.Range("K2").Value = 1178.2
$ws.Range("M2").Value = -1065.2
$ws.Range("H102").Value = 4378.1177
$ws.Range("I102").Value = 3173.3572
$ws.Range("K102").Value = 3173.3572
$ws.Range("M102").Value = -1551.3572
$ws.Range("H116").Value = 1178.2
$ws.Range("I116").Value = 1178.2
$ws.Range("K116").Value = 1178.2
$ws.Range("M116").Value = 1115.8
$ws.Range("H132").Value = 2395.889
$ws.Range("I132").Value = 2408.2354
$ws.Range("K132").Value = 7224.706200000001
$ws.Range("M132").Value = -4694.706200000001
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1178.2
$ws.Range("I3").Value = 1178.2
$ws.Range("K3").Value = 1178.2
$ws.Range("M3").Value = -1064.2
$ws.Range("H99").Value = 1728.2858
$ws.Range("I99").Value = 1779.8
$ws.Range("J99").Value = 1599.5
$ws.Range("K99").Value = 1779.8
$ws.Range("L99").Value = 1599.5
$ws.Range("M99").Value = -281.8
$ws.Range("N99").Value = -4595.5
$ws.Range("H134").Value = 3867.3076
$ws.Range("I134").Value = 3829
$ws.Range("K134").Value = 11487
$ws.Range("M134").Value = -8952
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4099.912
$ws.Range("I31").Value = 2020.5
$ws.Range("J31").Value = 4911.39
$ws.Range("K31").Value = 2020.5
$ws.Range("L31").Value = 4911.39
$ws.Range("M31").Value = -1725.5
$ws.Range("N31").Value = -5501.39
$ws.Range("H34").Value = 4099.912
$ws.Range("I34").Value = 2020.5
$ws.Range("J34").Value = 4911.39
$ws.Range("K34").Value = 2020.5
$ws.Range("L34").Value = 4911.39
$ws.Range("M34").Value = -1818.5
$ws.Range("N34").Value = -5315.39
$ws.Range("H132").Value = 3000
$ws.Range("I132").Value = 3000
$ws.Range("K132").Value = 9000
$ws.Range("M132").Value = -6470
$ws.Range("H134").Value = 454.22726
$ws.Range("I134").Value = 428.2381
$ws.Range("K134").Value = 1284.7143
$ws.Range("M134").Value = 1250.2857
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 3561.75
$ws.Range("I50").Value = 3482.3333
$ws.Range("J50").Value = 3800
$ws.Range("K50").Value = 10446.9999
$ws.Range("L50").Value = 11400
$ws.Range("M50").Value = -9965.999899999999
$ws.Range("N50").Value = -12362
$ws.Range("H53").Value = 3561.75
$ws.Range("I53").Value = 3482.3333
$ws.Range("J53").Value = 3800
$ws.Range("K53").Value = 10446.9999
$ws.Range("L53").Value = 11400
$ws.Range("M53").Value = -9965.999899999999
$ws.Range("N53").Value = -12362
$ws.Range("H68").Value = 3166.6667
$ws.Range("I68").Value = 1499
$ws.Range("K68").Value = 4497
$ws.Range("M68").Value = -3686
$ws.Range("H71").Value = 3166.6667
$ws.Range("I71").Value = 1499
$ws.Range("K71").Value = 13491
$ws.Range("M71").Value = -9435
$ws.Range("H128").Value = 251444.25
$ws.Range("I128").Value = 251444.25
$ws.Range("K128").Value = 754332.75
$ws.Range("M128").Value = -749352.75
$ws.Range("H131").Value = 1709.2972
$ws.Range("J131").Value = 2315.8333
$ws.Range("L131").Value = 6947.499899999999
$ws.Range("N131").Value = -17027.4999
$ws.Range("H136").Value = 7558.3335
$ws.Range("I136").Value = 6975
$ws.Range("K136").Value = 20925
$ws.Range("M136").Value = -15825
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H80").Value = 6183.8887
$ws.Range("I80").Value = 6231
$ws.Range("J80").Value = 6125
$ws.Range("K80").Value = 6231
$ws.Range("L80").Value = 6125
$ws.Range("M80").Value = -5233
$ws.Range("N80").Value = -8121
$ws.Range("H83").Value = 6183.8887
$ws.Range("I83").Value = 6231
$ws.Range("J83").Value = 6125
$ws.Range("K83").Value = 31155
$ws.Range("L83").Value = 30625
$ws.Range("M83").Value = -26163
$ws.Range("N83").Value = -40609
$ws.Range("H122").Value = 10000
$ws.Range("I122").Value = 10000
$ws.Range("K122").Value = 30000
$ws.Range("M122").Value = -27550
$ws.Range("H132").Value = 3657.4
$ws.Range("I132").Value = 3657.4
$ws.Range("K132").Value = 10972.2
$ws.Range("M132").Value = -8442.200000000001
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 712.55554
$ws.Range("J16").Value = 473.66666
$ws.Range("L16").Value = 473.66666
$ws.Range("N16").Value = -813.66666
$ws.Range("H61").Value = 3605.8333
$ws.Range("I61").Value = 2727
$ws.Range("J61").Value = 8000
$ws.Range("K61").Value = 2727
$ws.Range("L61").Value = 8000
$ws.Range("M61").Value = -2525
$ws.Range("N61").Value = -8404
$ws.Range("H93").Value = 6236.75
$ws.Range("I93").Value = 4900
$ws.Range("K93").Value = 4900
$ws.Range("M93").Value = -3652
$ws.Range("H100").Value = 8057.6924
$ws.Range("I100").Value = 5792.1665
$ws.Range("K100").Value = 5792.1665
$ws.Range("M100").Value = -5251.1665
$ws.Range("H113").Value = 3605.8333
$ws.Range("I113").Value = 2727
$ws.Range("J113").Value = 8000
$ws.Range("K113").Value = 2727
$ws.Range("L113").Value = 8000
$ws.Range("M113").Value = -557
$ws.Range("N113").Value = -12340
$ws.Range("H132").Value = 18989.25
$ws.Range("J132").Value = 19329
$ws.Range("L132").Value = 57987
$ws.Range("N132").Value = -63047
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 29313
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 29313
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 29313
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -30295
$ws.Range("H122").Value = 3209.75
$ws.Range("I122").Value = 3382.5715
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 10147.7145
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -7697.7145
$ws.Range("N122").Value = -10900
$ws.Range("H132").Value = 1853.4286
$ws.Range("I132").Value = 1853.4286
$ws.Range("K132").Value = 5560.2858
$ws.Range("M132").Value = -3030.2858
